# refactor: change namespace to Tamarin; bug fix: icon
#
# The "broke the icon again" task (Id 47) is completed: remove it from the
# Active list and re-add it to the top of the Inactive list, now marked
# Done with a completion date.

$wb = $excel.ActiveWorkbook

$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")

# Remove the completed task from the Active sheet (row 2: Id 47 "broke the
# icon again"). Everything below shifts up one row.
$active.Rows.Item(2).Delete()

# Insert a new row at the top of the Inactive sheet's data and record the
# same task there, now Done.
$inactive.Rows.Item(2).Insert()

$inactive.Cells.Item(2, 1).Value = 47
$inactive.Cells.Item(2, 2).Value = "broke the icon again"
$inactive.Cells.Item(2, 3).Value = "Done"
$inactive.Cells.Item(2, 4).Value = "Bug"
$inactive.Cells.Item(2, 5).Value = "'9/12/2018"
$inactive.Cells.Item(2, 6).Value = "'9/14/2018"

# Normalize the new row back to the workbook's plain default style (the
# leading apostrophes above keep the date-like text from being reinterpreted
# as real dates while still letting us drop back to the default format).
$inactive.Range("A2:F2").Style = "Normal"
